# Finish implementation of item transfer node manager:
# adds a new "Artisan" worksheet (hotkey / build reference sheet) and
# nudges the view-state (selection / scroll position) of a few other sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update view-state (selection) on the existing sheets, in the same
#    order a user tabbing through the workbook would touch them. Doing
#    this before the new sheet is added/activated means only the very
#    last sheet we activate ends up "tabSelected".
# ---------------------------------------------------------------------

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("C34").Select()

$ws4tab = $wb.Worksheets.Item("Sheet4")
$ws4tab.Activate()
$ws4tab.Range("D32").Select()

$wsUnitStats = $wb.Worksheets.Item("Unit stats")
$wsUnitStats.Activate()
$wsUnitStats.Range("P28").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Artisan" worksheet as the last tab.
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Artisan"

$ws.Columns.Item(3).ColumnWidth = 18.166666666666668

# Title
$ws.Range("A1").Value = "Artisan"

# Hotkey reference block 1 (Q/W/E/R - A - S/D/F - Z/X/C) with building names
$ws.Range("A3").Value = "Q"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Crude Axe"

$ws.Range("A4").Value = "W"
$ws.Range("C4").Value = "Crude Pickaxe"

$ws.Range("A5").Value = "E"
$ws.Range("C5").Value = "Mineshaft"

$ws.Range("A6").Value = "R"
$ws.Range("C6").Value = "Hell Forge"

$ws.Range("A8").Value = "A"
$ws.Range("C8").Value = "Workstation"

$ws.Range("A9").Value = "S"
$ws.Range("C9").Value = "Assembler"

$ws.Range("A10").Value = "D"
$ws.Range("C10").Value = "Minecart"

$ws.Range("A11").Value = "F"

$ws.Range("A13").Value = "Z"
$ws.Range("C13").Value = "Furnace"

$ws.Range("A14").Value = "X"
$ws.Range("C14").Value = "Furnace"

$ws.Range("A15").Value = "C"

# Hotkey reference block 2 - "Prospector" column
$ws.Range("F1").Value = "Prospector"

$ws.Range("F3").Value = "Q"
$ws.Range("H3").Value = "Defile"

$ws.Range("F4").Value = "W"
$ws.Range("H4").Value = "Extract Fel"

$ws.Range("F5").Value = "E"
$ws.Range("H5").Value = "Demonfruit"

$ws.Range("F6").Value = "R"
$ws.Range("H6").Value = "Fel Basin"

$ws.Range("F8").Value = "A"

$ws.Range("F9").Value = "S"
$ws.Range("H9").Value = "Scouts"

$ws.Range("F10").Value = "D"
$ws.Range("H10").Value = "Transfer Fel"

$ws.Range("F11").Value = "F"

$ws.Range("F13").Value = "Z"
$ws.Range("F14").Value = "X"
$ws.Range("F15").Value = "C"

# Hotkey reference block 3 - "Researcher" column
$ws.Range("K1").Value = "Researcher"

$ws.Range("K3").Value = "Q"
$ws.Range("K4").Value = "W"
$ws.Range("K5").Value = "E"
$ws.Range("K6").Value = "R"
$ws.Range("K8").Value = "A"
$ws.Range("K9").Value = "S"
$ws.Range("K10").Value = "D"
$ws.Range("K11").Value = "F"
$ws.Range("K13").Value = "Z"
$ws.Range("K14").Value = "X"
$ws.Range("K15").Value = "C"

# Research section
$ws.Range("A18").Value = "Research"

$ws.Range("A19").Value = "Q"
$ws.Range("C19").Value = "Tanks"
$ws.Range("D19").Value = "4 Iron, 4 Iron, 3 Iron + Frame"

$ws.Range("A20").Value = "W"
$ws.Range("C20").Value = "Converters"

$ws.Range("A21").Value = "E"
$ws.Range("C21").Value = "Automaton"
$ws.Range("D21").Value = "2x Iron, 2x Log, Rabbit + Frog"

$ws.Range("A22").Value = "R"
$ws.Range("C22").Value = "Depot"
$ws.Range("D22").Value = "Frame + Iron, Copper + Frog,  Skink + Rabbit"

$ws.Range("A24").Value = "A"
$ws.Range("A25").Value = "S"
$ws.Range("A26").Value = "D"
$ws.Range("A27").Value = "F"

$ws.Range("A29").Value = "Z"
$ws.Range("A30").Value = "X"
$ws.Range("A31").Value = "C"

# Automation notes
$ws.Range("B33").Value = "To Automate:"

$ws.Range("C34").Value = "Stone mining"
$ws.Range("E34").Value = "Quarry"

$ws.Range("C35").Value = "Log mining"
$ws.Range("E35").Value = "Mutator"

$ws.Range("C36").Value = "Steel processing"
$ws.Range("E36").Value = "Foundry"

# Merge/left-align the title
$ws.Range("A1:B1").HorizontalAlignment = -4131
$ws.Range("A1:B1").MergeCells = $true

# Make the new sheet the active tab / selection, like a freshly-edited sheet
$ws.Activate()
$ws.Range("R8").Select()
